# Applies:
#  1. Date placeholder field text "05-03-2024" -> "17-06-2024" on every
#     slide layout and the slide master (the cached text of the
#     datetimeFigureOut footer field).
#  2. "Analyzer" -> "Playback" in the second slide's logo lockup text.

$p = $ppt.ActivePresentation

function Update-DatePlaceholder($shapes) {
    for ($i = 1; $i -le $shapes.Count; $i++) {
        $sh = $shapes.Item($i)
        if ($sh.HasTextFrame -and $sh.TextFrame.HasText) {
            $tr = $sh.TextFrame.TextRange
            if ($tr.Text -eq "05-03-2024") {
                $tr.Text = "17-06-2024"
            }
        }
    }
}

# Slide master date placeholder.
Update-DatePlaceholder $p.SlideMaster.Shapes

# Every slide layout's date placeholder.
$layouts = $p.SlideMaster.CustomLayouts
for ($j = 1; $j -le $layouts.Count; $j++) {
    Update-DatePlaceholder $layouts.Item($j).Shapes
}

# Slide 2: "Analyzer" -> "Playback" inside the grouped logo lockup.
$s = $p.Slides.Item(2)
$grp = $s.Shapes.Item(2)
$rect = $grp.GroupItems.Item(1)
$tr = $rect.TextFrame.TextRange
$para = $tr.Paragraphs(2)
if ($para.Text -eq "Analyzer") {
    $para.Text = "Playback"
}
